$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196, shifting existing rows 196:205 down to 197:206
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row 196 with the new weekly record
$ws.Cells.Item(196, 1).Value = 11
$ws.Cells.Item(196, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(196, 3).Value = "Bíobío"
$ws.Cells.Item(196, 4).Value = 44747
$ws.Cells.Item(196, 5).Value = 8
$ws.Cells.Item(196, 6).Value = 100112040
$ws.Cells.Item(196, 7).Value = "Cilantro"
$ws.Cells.Item(196, 8).Value = "Sin especificar"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 70
$ws.Cells.Item(196, 11).Value = 9000
$ws.Cells.Item(196, 12).Value = 10000
$ws.Cells.Item(196, 13).Value = 9286
$ws.Cells.Item(196, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(196, 15).Value = "Región Metropolitana"
$ws.Cells.Item(196, 16).Value = 258
$ws.Cells.Item(196, 17).Value = 36
$ws.Cells.Item(196, 18).Value = "Hortaliza"
